# Added new profile search test cases.
# Appends 4 new rows (FindProfileWithLastNameTest, FindProfileWithRoleTest,
# FindProfileWithPrimaryInstitutionTest, FindProfileWithCountryTest) to the
# "Test Cases" sheet, right after the existing "ProfileFollowerTest" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 4 (ViewProfileTest) already carries the exact border/fill/font
# combination used by the plain (non-header, non-wrap) data rows, so copy
# its formatting down onto the four new rows before filling in the values.
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A14:D17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$newRows = @(
    @("FindProfileWithLastNameTest", "Verify that user is able to find other profiles with their last name", "Y", "SKIP"),
    @("FindProfileWithRoleTest", "Verify that user is able to find other profiles with their Title/Role", "Y", "SKIP"),
    @("FindProfileWithPrimaryInstitutionTest", "Verify that user is able to find other profiles with their Primary Institution", "Y", "PASS"),
    @("FindProfileWithCountryTest", "Verify that user is able to find other profiles with their Country", "Y", "SKIP")
)

$r = 14
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws.Range("C17").Select() | Out-Null
